# Update a handful of imputed values in the RandomForest result sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value  = -21.11390000000001
$ws.Range("A10").Value = -20.48889999999997
$ws.Range("A12").Value = -22.55530000000003
$ws.Range("C13").Value = -12.58409999999999
$ws.Range("A18").Value = -22.44290000000003
